$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 17:06"

# Update country data rows (name + Casos totales / Nuevos casos / Casos activos / Recuperados / Casos criticos / Muertes hoy / Muertes).
# This reflects the refreshed COVID-19 snapshot; a few countries changed rank order, so the country
# label in column A is re-set alongside the numbers to keep everything aligned.

$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 8228001
$ws.Range("C4").Value = 11686
$ws.Range("D4").Value = 5328605
$ws.Range("E4").Value = 2676461
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 218
$ws.Range("H4").Value = 222935

$ws.Range("A5").Value = "India"
$ws.Range("B5").Value = 7383104
$ws.Range("C5").Value = 17595
$ws.Range("D5").Value = 6466336
$ws.Range("E5").Value = 804411
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 211
$ws.Range("H5").Value = 112357

$ws.Range("A17").Value = "Chile"
$ws.Range("B17").Value = 488190
$ws.Range("C17").Value = 1694
$ws.Range("D17").Value = 461097
$ws.Range("E17").Value = 13564
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 95
$ws.Range("H17").Value = 13529

$ws.Range("A21").Value = "Alemania"
$ws.Range("B21").Value = 353822
$ws.Range("C21").Value = 5006
$ws.Range("D21").Value = 284600
$ws.Range("E21").Value = 59393
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 19
$ws.Range("H21").Value = 9829

$ws.Range("A22").Value = "Indonesia"
$ws.Range("B22").Value = 353461
$ws.Range("C22").Value = 4301
$ws.Range("D22").Value = 277544
$ws.Range("E22").Value = 63570
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 79
$ws.Range("H22").Value = 12347

$ws.Range("A30").Value = "Canada"
$ws.Range("B30").Value = 192442
$ws.Range("C30").Value = 710
$ws.Range("D30").Value = 162203
$ws.Range("E30").Value = 20531
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 9
$ws.Range("H30").Value = 9708

$ws.Range("A31").Value = "Belgica"
$ws.Range("B31").Value = 191959
$ws.Range("C31").Value = 10448
$ws.Range("D31").Value = 20720
$ws.Range("E31").Value = 160912
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 49
$ws.Range("H31").Value = 10327

$ws.Range("A48").Value = "Guatemala"
$ws.Range("B48").Value = 100431
$ws.Range("C48").Value = 666
$ws.Range("D48").Value = 89494
$ws.Range("E48").Value = 7459
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 25
$ws.Range("H48").Value = 3478

$ws.Range("A51").Value = "Japon"
$ws.Range("B51").Value = 91431
$ws.Range("C51").Value = 721
$ws.Range("D51").Value = 84451
$ws.Range("E51").Value = 5330
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 4
$ws.Range("H51").Value = 1650

$ws.Range("A59").Value = "Moldavia"
$ws.Range("B59").Value = 65860
$ws.Range("C59").Value = 784
$ws.Range("D59").Value = 46543
$ws.Range("E59").Value = 17768
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 19
$ws.Range("H59").Value = 1549

$ws.Range("A65").Value = "Singapur"
$ws.Range("B65").Value = 57901
$ws.Range("C65").Value = 9
$ws.Range("D65").Value = 57784
$ws.Range("E65").Value = 89
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 28

$ws.Range("A74").Value = "Kenia"
$ws.Range("B74").Value = 43580
$ws.Range("C74").Value = 437
$ws.Range("D74").Value = 31648
$ws.Range("E74").Value = 11119
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 8
$ws.Range("H74").Value = 813

$ws.Range("A80").Value = "Birmania"
$ws.Range("B80").Value = 33488
$ws.Range("C80").Value = 1137
$ws.Range("D80").Value = 15477
$ws.Range("E80").Value = 17212
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 34
$ws.Range("H80").Value = 799

$ws.Range("A81").Value = "Jordania"
$ws.Range("B81").Value = 33009
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 6565
$ws.Range("E81").Value = 26162
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 282

$ws.Range("A82").Value = "Bosnia y Herzegovina"
$ws.Range("B82").Value = 32845
$ws.Range("C82").Value = 621
$ws.Range("D82").Value = 24603
$ws.Range("E82").Value = 7262
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 8
$ws.Range("H82").Value = 980

$ws.Range("A95").Value = "Albania"
$ws.Range("B95").Value = 16501
$ws.Range("C95").Value = 289
$ws.Range("D95").Value = 9957
$ws.Range("E95").Value = 6101
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 4
$ws.Range("H95").Value = 443

$ws.Range("A115").Value = "Jamaica"
$ws.Range("B115").Value = 8132
$ws.Range("C115").Value = 65
$ws.Range("D115").Value = 3653
$ws.Range("E115").Value = 4317
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 162

$ws.Range("A116").Value = "Zimbabue"
$ws.Range("B116").Value = 8075
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 7669
$ws.Range("E116").Value = 175
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 231

$ws.Range("A122").Value = "Cuba"
$ws.Range("B122").Value = 6118
$ws.Range("C122").Value = 56
$ws.Range("D122").Value = 5702
$ws.Range("E122").Value = 292
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 124

$ws.Range("A127").Value = "Sri Lanka"
$ws.Range("B127").Value = 5354
$ws.Range("C127").Value = 110
$ws.Range("D127").Value = 3385
$ws.Range("E127").Value = 1956
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 13

$ws.Range("A128").Value = "Nicaragua"
$ws.Range("B128").Value = 5353
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 4225
$ws.Range("E128").Value = 974
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 154
